# Update the nowcast table with the latest run: revise existing vintages
# (rows 2-11) and append the newest vintage (row 12), per "add results from
# latest run".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:K1) already holds the correct labels ("Row", "Prognose",
# ... "Revision") and is left untouched.

# Column A rows 2-11 already hold the correct vintage-date text
# ("2025-03-30" .. "2025-08-15") and are left untouched. Row 12 is the
# newly added vintage from the latest run; force text format first so
# the date-looking string is stored as text (matching the rest of the
# column) instead of being auto-converted to an Excel date serial.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '2025-08-30'

# Data matrix B2:K12 (nowcast value + per-block revision contributions)
$ws.Range("B2").Value = 0.32302385808069306
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("B3").Value = 0.30919947782404233
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.0016898276287108719
$ws.Range("E3").Value = 0.00014075689805993782
$ws.Range("F3").Value = -0.0007291553961898643
$ws.Range("G3").Value = 0.00039154866702117915
$ws.Range("H3").Value = 0.000022095875308363403
$ws.Range("I3").Value = -0.0011025162175083874
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.004579355338180824
$ws.Range("B4").Value = 0.30781136843838536
$ws.Range("C4").Value = -0.0006780723558998704
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.00004293300605324955
$ws.Range("F4").Value = -0.000004332504420695083
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.000024693160090091088
$ws.Range("I4").Value = -0.000595145764247854
$ws.Range("J4").Value = -0.00005748034623903365
$ws.Range("K4").Value = -0.00041144236860679584
$ws.Range("B5").Value = 0.2821447458739408
$ws.Range("C5").Value = 0.0030370782775787776
$ws.Range("D5").Value = -0.003135596855176412
$ws.Range("E5").Value = -0.000060748013382638945
$ws.Range("F5").Value = -0.0009786365304629086
$ws.Range("G5").Value = -0.0009856342969447047
$ws.Range("H5").Value = -0.00010067756016386232
$ws.Range("I5").Value = -0.00020512749434252098
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00031224912210831013
$ws.Range("B6").Value = 0.36342496307109823
$ws.Range("C6").Value = 0.021100256859687014
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = -0.00010619380307271246
$ws.Range("F6").Value = -0.00006388232391242386
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.00005492604009847419
$ws.Range("I6").Value = -0.0017777353708719739
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = -0.000509612226942513
$ws.Range("B7").Value = 0.3614392734724381
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.0020075126454683285
$ws.Range("E7").Value = -0.0001898219393004642
$ws.Range("F7").Value = -0.001212035563656326
$ws.Range("G7").Value = 0.0006950348481490643
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.00009382099003910959
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.0016150895666742215
$ws.Range("B8").Value = 0.2659459813734798
$ws.Range("C8").Value = -0.023519605164103778
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.00005492713531008989
$ws.Range("F8").Value = -0.00040975135487218906
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.000014420778755051014
$ws.Range("I8").Value = 0.0012377999295085905
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.0005469371480167351
$ws.Range("B9").Value = 0.2601711448942538
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.0012966781285655672
$ws.Range("E9").Value = -0.0023147939383170823
$ws.Range("F9").Value = -0.006283944970567324
$ws.Range("G9").Value = 0.0006166702378112915
$ws.Range("H9").Value = -0.00019736592570409703
$ws.Range("I9").Value = -0.0002385430470102691
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.0004490245641920465
$ws.Range("B10").Value = 0.4005317629620292
$ws.Range("C10").Value = 0.0474083528187003
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.000026171642466790374
$ws.Range("F10").Value = -0.00039715574072223084
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.000026902045027001622
$ws.Range("I10").Value = -0.0004538251416684124
$ws.Range("J10").Value = -0.0023788044561584144
$ws.Range("K10").Value = 0.000909258556880399
$ws.Range("B11").Value = 0.37636991322927116
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -0.0017306195834714908
$ws.Range("E11").Value = 0.0006432450284716704
$ws.Range("F11").Value = -0.00031537052722860744
$ws.Range("G11").Value = 0.0014527552225283926
$ws.Range("H11").Value = 0.000135512287960789
$ws.Range("I11").Value = -0.004724537595136337
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0.013843182496157569
$ws.Range("B12").Value = 0.3161818873807198
$ws.Range("C12").Value = -0.030710080589056766
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.00008103142360267607
$ws.Range("F12").Value = 0.00001233246134499027
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.000005545158359644843
$ws.Range("I12").Value = -0.0006975689218262495
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.0033268278891097625
